$wb = $excel.ActiveWorkbook

# ---- 1. Rename headers on existing sheets ----
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# ---- 2. Add the new 'PO Forecast' sheet after 'Monthly Trend' ----
$newSheet = $wb.Worksheets.Add($null, $wsMonthly)
$newSheet.Name = "PO Forecast"

# Match page margins used by the other sheets (inches -> points: 1in = 72pt)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---- 3. Copy the header / date-column formatting from 'Weekly Quantity' ----
# Header style (bold, centered, bordered) -> row 1, columns A:D
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$newSheet.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# Date-column style (custom date/time number format) -> column A, rows 2:48
# (source range is the populated A2:A40 on 'Weekly Quantity'; PasteSpecial tiles it
#  to cover the larger A2:A48 destination)
$wsWeekly.Range("A2:A40").Copy() | Out-Null
$newSheet.Range("A2:A48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- 4. Header values ----
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# ---- 5. Forecast data rows ----
$newSheet.Cells.Item(2, 1).Value = 45249.99999999999
$newSheet.Cells.Item(2, 2).Value = 387
$newSheet.Cells.Item(2, 3).Value = -500.5026098737151
$newSheet.Cells.Item(2, 4).Value = 1239.282428803995
$newSheet.Cells.Item(3, 1).Value = 45256.99999999999
$newSheet.Cells.Item(3, 2).Value = 400
$newSheet.Cells.Item(3, 3).Value = -557.2002479944309
$newSheet.Cells.Item(3, 4).Value = 1237.624755297643
$newSheet.Cells.Item(4, 1).Value = 45270.99999999999
$newSheet.Cells.Item(4, 2).Value = 427
$newSheet.Cells.Item(4, 3).Value = -471.6440586408746
$newSheet.Cells.Item(4, 4).Value = 1270.24434384777
$newSheet.Cells.Item(5, 1).Value = 45277.99999999999
$newSheet.Cells.Item(5, 2).Value = 440
$newSheet.Cells.Item(5, 3).Value = -375.7989253505222
$newSheet.Cells.Item(5, 4).Value = 1303.471919743216
$newSheet.Cells.Item(6, 1).Value = 45298.99999999999
$newSheet.Cells.Item(6, 2).Value = 479
$newSheet.Cells.Item(6, 3).Value = -407.9036826363477
$newSheet.Cells.Item(6, 4).Value = 1411.905712508739
$newSheet.Cells.Item(7, 1).Value = 45305.99999999999
$newSheet.Cells.Item(7, 2).Value = 492
$newSheet.Cells.Item(7, 3).Value = -330.5786163673868
$newSheet.Cells.Item(7, 4).Value = 1424.474897071015
$newSheet.Cells.Item(8, 1).Value = 45312.99999999999
$newSheet.Cells.Item(8, 2).Value = 506
$newSheet.Cells.Item(8, 3).Value = -404.7015985903642
$newSheet.Cells.Item(8, 4).Value = 1394.863977295038
$newSheet.Cells.Item(9, 1).Value = 45319.99999999999
$newSheet.Cells.Item(9, 2).Value = 519
$newSheet.Cells.Item(9, 3).Value = -328.6306812883897
$newSheet.Cells.Item(9, 4).Value = 1459.627891359189
$newSheet.Cells.Item(10, 1).Value = 45326.99999999999
$newSheet.Cells.Item(10, 2).Value = 532
$newSheet.Cells.Item(10, 3).Value = -341.3395572138847
$newSheet.Cells.Item(10, 4).Value = 1439.32499523211
$newSheet.Cells.Item(11, 1).Value = 45333.99999999999
$newSheet.Cells.Item(11, 2).Value = 545
$newSheet.Cells.Item(11, 3).Value = -369.9401990944316
$newSheet.Cells.Item(11, 4).Value = 1408.73490488371
$newSheet.Cells.Item(12, 1).Value = 45340.99999999999
$newSheet.Cells.Item(12, 2).Value = 558
$newSheet.Cells.Item(12, 3).Value = -326.9983743275707
$newSheet.Cells.Item(12, 4).Value = 1447.85915075067
$newSheet.Cells.Item(13, 1).Value = 45347.99999999999
$newSheet.Cells.Item(13, 2).Value = 571
$newSheet.Cells.Item(13, 3).Value = -379.4857761904381
$newSheet.Cells.Item(13, 4).Value = 1403.510525024212
$newSheet.Cells.Item(14, 1).Value = 45354.99999999999
$newSheet.Cells.Item(14, 2).Value = 585
$newSheet.Cells.Item(14, 3).Value = -319.680470810392
$newSheet.Cells.Item(14, 4).Value = 1455.66345916463
$newSheet.Cells.Item(15, 1).Value = 45368.99999999999
$newSheet.Cells.Item(15, 2).Value = 611
$newSheet.Cells.Item(15, 3).Value = -290.8691800108464
$newSheet.Cells.Item(15, 4).Value = 1473.342204878539
$newSheet.Cells.Item(16, 1).Value = 45375.99999999999
$newSheet.Cells.Item(16, 2).Value = 624
$newSheet.Cells.Item(16, 3).Value = -260.9315583621005
$newSheet.Cells.Item(16, 4).Value = 1508.970792518178
$newSheet.Cells.Item(17, 1).Value = 45382.99999999999
$newSheet.Cells.Item(17, 2).Value = 637
$newSheet.Cells.Item(17, 3).Value = -279.8207052325578
$newSheet.Cells.Item(17, 4).Value = 1527.261228117488
$newSheet.Cells.Item(18, 1).Value = 45389.99999999999
$newSheet.Cells.Item(18, 2).Value = 650
$newSheet.Cells.Item(18, 3).Value = -270.9407193001948
$newSheet.Cells.Item(18, 4).Value = 1565.092142511081
$newSheet.Cells.Item(19, 1).Value = 45396.99999999999
$newSheet.Cells.Item(19, 2).Value = 664
$newSheet.Cells.Item(19, 3).Value = -213.9994309215312
$newSheet.Cells.Item(19, 4).Value = 1510.175274285855
$newSheet.Cells.Item(20, 1).Value = 45403.99999999999
$newSheet.Cells.Item(20, 2).Value = 677
$newSheet.Cells.Item(20, 3).Value = -194.2269456018034
$newSheet.Cells.Item(20, 4).Value = 1587.532561529947
$newSheet.Cells.Item(21, 1).Value = 45410.99999999999
$newSheet.Cells.Item(21, 2).Value = 690
$newSheet.Cells.Item(21, 3).Value = -165.3160398198987
$newSheet.Cells.Item(21, 4).Value = 1582.08392062195
$newSheet.Cells.Item(22, 1).Value = 45417.99999999999
$newSheet.Cells.Item(22, 2).Value = 703
$newSheet.Cells.Item(22, 3).Value = -216.5142408286855
$newSheet.Cells.Item(22, 4).Value = 1577.106733107938
$newSheet.Cells.Item(23, 1).Value = 45424.99999999999
$newSheet.Cells.Item(23, 2).Value = 716
$newSheet.Cells.Item(23, 3).Value = -131.2099173958073
$newSheet.Cells.Item(23, 4).Value = 1563.949207558399
$newSheet.Cells.Item(24, 1).Value = 45438.99999999999
$newSheet.Cells.Item(24, 2).Value = 743
$newSheet.Cells.Item(24, 3).Value = -111.6254775924856
$newSheet.Cells.Item(24, 4).Value = 1581.223081037909
$newSheet.Cells.Item(25, 1).Value = 45445.99999999999
$newSheet.Cells.Item(25, 2).Value = 756
$newSheet.Cells.Item(25, 3).Value = -103.2497400591338
$newSheet.Cells.Item(25, 4).Value = 1625.102153217308
$newSheet.Cells.Item(26, 1).Value = 45459.99999999999
$newSheet.Cells.Item(26, 2).Value = 782
$newSheet.Cells.Item(26, 3).Value = -110.0596987650745
$newSheet.Cells.Item(26, 4).Value = 1702.523268687125
$newSheet.Cells.Item(27, 1).Value = 45466.99999999999
$newSheet.Cells.Item(27, 2).Value = 795
$newSheet.Cells.Item(27, 3).Value = -115.5639702080831
$newSheet.Cells.Item(27, 4).Value = 1699.903712496686
$newSheet.Cells.Item(28, 1).Value = 45473.99999999999
$newSheet.Cells.Item(28, 2).Value = 808
$newSheet.Cells.Item(28, 3).Value = -91.35347417998462
$newSheet.Cells.Item(28, 4).Value = 1706.636747947174
$newSheet.Cells.Item(29, 1).Value = 45480.99999999999
$newSheet.Cells.Item(29, 2).Value = 821
$newSheet.Cells.Item(29, 3).Value = -95.41441774546628
$newSheet.Cells.Item(29, 4).Value = 1693.994378058694
$newSheet.Cells.Item(30, 1).Value = 45487.99999999999
$newSheet.Cells.Item(30, 2).Value = 835
$newSheet.Cells.Item(30, 3).Value = -25.20227891059661
$newSheet.Cells.Item(30, 4).Value = 1725.063198668305
$newSheet.Cells.Item(31, 1).Value = 45508.99999999999
$newSheet.Cells.Item(31, 2).Value = 874
$newSheet.Cells.Item(31, 3).Value = -28.71299767614274
$newSheet.Cells.Item(31, 4).Value = 1751.204165210429
$newSheet.Cells.Item(32, 1).Value = 45529.99999999999
$newSheet.Cells.Item(32, 2).Value = 914
$newSheet.Cells.Item(32, 3).Value = 0.7086782529100857
$newSheet.Cells.Item(32, 4).Value = 1831.141176941469
$newSheet.Cells.Item(33, 1).Value = 45550.99999999999
$newSheet.Cells.Item(33, 2).Value = 953
$newSheet.Cells.Item(33, 3).Value = 84.28853452290841
$newSheet.Cells.Item(33, 4).Value = 1845.804972597981
$newSheet.Cells.Item(34, 1).Value = 45557.99999999999
$newSheet.Cells.Item(34, 2).Value = 966
$newSheet.Cells.Item(34, 3).Value = 110.7073745482578
$newSheet.Cells.Item(34, 4).Value = 1794.32451085003
$newSheet.Cells.Item(35, 1).Value = 45571.99999999999
$newSheet.Cells.Item(35, 2).Value = 993
$newSheet.Cells.Item(35, 3).Value = 85.64902140291684
$newSheet.Cells.Item(35, 4).Value = 1877.88872482388
$newSheet.Cells.Item(36, 1).Value = 45578.99999999999
$newSheet.Cells.Item(36, 2).Value = 1006
$newSheet.Cells.Item(36, 3).Value = 155.9597905173737
$newSheet.Cells.Item(36, 4).Value = 1930.711849745087
$newSheet.Cells.Item(37, 1).Value = 45585.99999999999
$newSheet.Cells.Item(37, 2).Value = 1019
$newSheet.Cells.Item(37, 3).Value = 87.01921637679651
$newSheet.Cells.Item(37, 4).Value = 1919.531609941786
$newSheet.Cells.Item(38, 1).Value = 45592.99999999999
$newSheet.Cells.Item(38, 2).Value = 1032
$newSheet.Cells.Item(38, 3).Value = 192.0061666987449
$newSheet.Cells.Item(38, 4).Value = 1974.746961991479
$newSheet.Cells.Item(39, 1).Value = 45599.99999999999
$newSheet.Cells.Item(39, 2).Value = 1045
$newSheet.Cells.Item(39, 3).Value = 101.8450631541447
$newSheet.Cells.Item(39, 4).Value = 1939.021239326717
$newSheet.Cells.Item(40, 1).Value = 45620.99999999999
$newSheet.Cells.Item(40, 2).Value = 1085
$newSheet.Cells.Item(40, 3).Value = 197.8268139584969
$newSheet.Cells.Item(40, 4).Value = 1981.222432550885
$newSheet.Cells.Item(41, 1).Value = 45627.99999999999
$newSheet.Cells.Item(41, 2).Value = 1098
$newSheet.Cells.Item(41, 3).Value = 224.5935073215157
$newSheet.Cells.Item(41, 4).Value = 1967.351482699546
$newSheet.Cells.Item(42, 1).Value = 45634.99999999999
$newSheet.Cells.Item(42, 2).Value = 1111
$newSheet.Cells.Item(42, 3).Value = 214.9211588938277
$newSheet.Cells.Item(42, 4).Value = 1994.298043053034
$newSheet.Cells.Item(43, 1).Value = 45641.99999999999
$newSheet.Cells.Item(43, 2).Value = 1124
$newSheet.Cells.Item(43, 3).Value = 254.5513916249247
$newSheet.Cells.Item(43, 4).Value = 2022.249545332515
$newSheet.Cells.Item(44, 1).Value = 45648.99999999999
$newSheet.Cells.Item(44, 2).Value = 1137
$newSheet.Cells.Item(44, 3).Value = 241.1089963303712
$newSheet.Cells.Item(44, 4).Value = 2008.717029683696
$newSheet.Cells.Item(45, 1).Value = 45655.99999999999
$newSheet.Cells.Item(45, 2).Value = 1150
$newSheet.Cells.Item(45, 3).Value = 263.6489984481734
$newSheet.Cells.Item(45, 4).Value = 2037.957188555484
$newSheet.Cells.Item(46, 1).Value = 45662.99999999999
$newSheet.Cells.Item(46, 2).Value = 1164
$newSheet.Cells.Item(46, 3).Value = 295.1424553250141
$newSheet.Cells.Item(46, 4).Value = 2053.397762616099
$newSheet.Cells.Item(47, 1).Value = 45669.99999999999
$newSheet.Cells.Item(47, 2).Value = 1177
$newSheet.Cells.Item(47, 3).Value = 314.3700880665366
$newSheet.Cells.Item(47, 4).Value = 2074.738923087642
$newSheet.Cells.Item(48, 1).Value = 45676.99999999999
$newSheet.Cells.Item(48, 2).Value = 1190
$newSheet.Cells.Item(48, 3).Value = 295.9272160100623
$newSheet.Cells.Item(48, 4).Value = 2072.927218528439
